# 2v2statsLeagueXL.xlsx - "Add files via upload" update
# Updates the per-player raw stat totals (K, HS, D, UD, EF, DAM, R) for the
# players in rows 8-15 of Sheet1. The dependent ratio formulas in columns
# I:O recalculate automatically from these raw values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ column letter = new value }
$updates = @{
    8  = @{ B = 44; C = 29; D = 36;            F = 17; G = 5040; H = 62 }
    9  = @{ B = 41; C = 19; D = 34; E = 107;    F = 14; G = 4866; H = 62 }
    10 = @{ B = 46; C = 14; D = 40; E = 176;            G = 5305; H = 73 }
    11 = @{ B = 47; C = 24; D = 46; E = 67;     F = 7;  G = 5416; H = 73 }
    12 = @{ B = 56; C = 38; D = 50; E = 400;    F = 29; G = 6996; H = 73 }
    13 = @{ B = 18; C = 9;  D = 53;             F = 16; G = 2380; H = 73 }
    14 = @{ B = 11; C = 3;  D = 48; E = 110;    F = 18; G = 2514; H = 55 }
    15 = @{ B = 21; C = 6;  D = 47;                     G = 2583; H = 55 }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    foreach ($col in $cells.Keys) {
        $ws.Range("$col$row").Value = $cells[$col]
    }
}

# Restore the window size/position metadata recorded when the file was
# last saved.
$excel.ActiveWindow.Left = 4935
$excel.ActiveWindow.Top = 5340

# Update the selected cell on Sheet1.
$ws.Range("K19").Select()
